$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:P1").UnMerge()

$title = $ws.Range("A1").Value
$ws.Range("F1").Value = $title
$ws.Range("A1").Clear()

# Row1 styling: B1:E1 and G1:P1 keep font1 18pt bold wrapText only (no horizontal center)
$ws.Range("B1:E1").HorizontalAlignment = 1
$ws.Range("B1:E1").WrapText = $true
$ws.Range("G1:P1").HorizontalAlignment = 1
$ws.Range("G1:P1").WrapText = $true

# F1: font1 18pt bold, no alignment at all (no wrap), clear center
$ws.Range("F1").HorizontalAlignment = 1
$ws.Range("F1").WrapText = $false

# Row2 changes
$ws.Range("A2:P2").UnMerge()
$dept = $ws.Range("A2").Value
$ws.Range("G2").Value = $dept
$ws.Range("A2").Clear()
$ws.Range("H2").Clear()
$ws.Range("L2:O2").Clear()

$ws.Range("B2:F2").HorizontalAlignment = 1
$ws.Range("G2").HorizontalAlignment = 1
$ws.Range("I2:K2").HorizontalAlignment = 1
$ws.Range("P2").HorizontalAlignment = 1

# Row3 changes: drop horizontal=center, keep text, font2 20pt (no alignment)
$ws.Range("A3:P3").UnMerge()
$ws.Range("A3:P3").HorizontalAlignment = 1
